$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 blank rows before row 18 to make room for the new "find/replace"
# related feature rows; this pushes the existing rows 18-49 down to 26-57.
$ws.Rows("18:25").Insert()

# Update the "语法高亮" entry to reflect that it is a simplified implementation.
$ws.Range("B14").Value = "(简单的)语法高亮"

# Fill in the newly inserted rows with the new feature rows.
$ws.Range("B18").Value = "分割布局支持"
$ws.Range("B19").Value = "括号匹配显示"
$ws.Range("B20").Value = "选中词高亮"
$ws.Range("B21").Value = "自动缩进"
$ws.Range("B22").Value = "符号自动补全/覆盖"
$ws.Range("A23").Value = "查找替换"
$ws.Range("B23").Value = "单文件查找/替换"
$ws.Range("B24").Value = "多文件查找"
$ws.Range("B25").Value = "正则表达式支持"

# The old "外部文件编辑变化监控" row (originally row 47) is now row 55 after
# the insert above; remove it entirely.
$ws.Rows("55:55").Delete()

# Update the view: scroll/selection now rests on B23, with no pinned top-left cell.
$ws.Range("B23").Select()
